# Loan RBI, Variable Instalments
# On the "Repayment schedule" sheet a new (blank) column is inserted
# before the old column N ("Late"), pushing the old N/O/P columns
# ("Late" / "heading" / "Outstanding") one position to the right
# (O/P/Q). The sheet is then left as the active sheet/tab with cell
# Q12 selected (previously "Edit Repayment Schedule" was the active tab).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at position N (14), shifting existing data right.
$ws.Columns.Item(14).Insert()

# The newly inserted column keeps the width of its left neighbour (column M).
$ws.Columns.Item(14).ColumnWidth = 10.17

# Make "Repayment schedule" the active sheet/tab and select Q12 on it
# (this also clears the tabSelected flag previously on
# "Edit Repayment Schedule").
$ws.Activate() | Out-Null
$ws.Range("Q12").Select() | Out-Null
